# SE-341 Added handling of OD600 data.
#
# The "openbis-metadata" sheet's example Strain value is updated from the
# placeholder "strain1" to "foo" (matching the other example placeholders
# in that column), and the active sheet / selection state left by the
# author is restored (the "openbis-metadata" sheet becomes the active tab
# with cell C11 selected).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("openbis-metadata")

# Update the example Strain value placeholder.
$ws1.Range("B3").Value = "foo"

# Make "openbis-metadata" the active sheet with C11 selected (this also
# clears the previously active selection/tab on "openbis-data").
$ws1.Activate()
$ws1.Range("C11").Select() | Out-Null
